$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (cohort_year 2020, period_index 5): num_customers 12 -> 13,
# retention_rate recalculated as num_customers / cohort_size (13/107)
$ws.Range("C7").Value = 13
$ws.Range("E7").Value = 0.1214953271028037

# Row 22 (cohort_year 2025, period_index 0): num_customers and cohort_size 51 -> 54
$ws.Range("C22").Value = 54
$ws.Range("D22").Value = 54
